# Connected Office Test Data - mark all CRUD tests as passed (TRUE)
# on the "Test Results" sheet, for rows 2-24 (columns B:E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Flip the Create/Read/Update/Delete test result flags from FALSE to TRUE
# for every data row (rows 2 through 24, columns B through E).
$ws.Range("B2:E24").Value = $true

# Reflect the scrolled/selected view state recorded for this sheet.
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("E2:E24").Select() | Out-Null
